$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price (D) and Volume (E) columns for rows 2-51
# so that numeric-looking strings (prices, percentages) are preserved exactly
# as text, matching the workbook's original inline-string storage.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '308.48'
$ws.Range("E2").Value = '1.04%'

# Row 3
$ws.Range("D3").Value = '38.50'
$ws.Range("E3").Value = '7.65%'

# Row 4
$ws.Range("D4").Value = '5.096'
$ws.Range("E4").Value = '1.23%'

# Row 5
$ws.Range("E5").Value = '1.10%'

# Row 6
$ws.Range("D6").Value = '1.979'
$ws.Range("E6").Value = '6.15%'

# Row 7
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = '7.938'
$ws.Range("E7").Value = '1.93%'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '0.9274'
$ws.Range("E8").Value = '0.80%'

# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '0.1446'
$ws.Range("E9").Value = '13.76%'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1958'
$ws.Range("E10").Value = '2.34%'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.09079'
$ws.Range("E11").Value = '0.45%'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.03507'
$ws.Range("E12").Value = '1.07%'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09812'
$ws.Range("E13").Value = '-0.51%'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001409'
$ws.Range("E14").Value = '-0.49%'

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.006036'
$ws.Range("E15").Value = '-3.20%'

# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.660'
$ws.Range("E16").Value = '-4.18%'

# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '4.205'
$ws.Range("E17").Value = '1.54%'

# Row 18
$ws.Range("E18").Value = '2.70%'

# Row 19
$ws.Range("D19").Value = '0.3431'
$ws.Range("E19").Value = '0.35%'

# Row 20
$ws.Range("E20").Value = '-1.06%'

# Row 21
$ws.Range("D21").Value = '4.801'
$ws.Range("E21").Value = '-8.44%'

# Row 22
$ws.Range("D22").Value = '0.2456'
$ws.Range("E22").Value = '6.47%'

# Row 23
$ws.Range("D23").Value = '0.04418'
$ws.Range("E23").Value = '-0.40%'

# Row 24
$ws.Range("D24").Value = '0.001218'
$ws.Range("E24").Value = '-1.37%'

# Row 25
$ws.Range("E25").Value = '4.78%'

# Row 27
$ws.Range("D27").Value = '0.0001303'
$ws.Range("E27").Value = '4.12%'

# Row 39
$ws.Range("D39").Value = '0.02099'
$ws.Range("E39").Value = '8.22%'

# Row 40
$ws.Range("D40").Value = '0.05150'
$ws.Range("E40").Value = '-2.57%'

# Row 41
$ws.Range("D41").Value = '0.007480'
$ws.Range("E41").Value = '-1.72%'

# Row 42
$ws.Range("D42").Value = '0.01015'
$ws.Range("E42").Value = '0.11%'

# Row 43
$ws.Range("D43").Value = '0.1357'
$ws.Range("E43").Value = '0.26%'

# Row 44
$ws.Range("D44").Value = '0.002145'
$ws.Range("E44").Value = '-0.80%'

# Row 45
$ws.Range("D45").Value = '0.009214'
$ws.Range("E45").Value = '-4.34%'

# Row 46
$ws.Range("D46").Value = '0.00006308'
$ws.Range("E46").Value = '2.80%'

# Row 47
$ws.Range("D47").Value = '0.00000000752'
$ws.Range("E47").Value = '0.12%'

# Row 48
$ws.Range("D48").Value = '0.003054'

# Row 49
$ws.Range("D49").Value = '0.001603'
$ws.Range("E49").Value = '-3.44%'

# Row 50
$ws.Range("D50").Value = '0.00002104'
$ws.Range("E50").Value = '0.12%'

# Row 51
$ws.Range("D51").Value = '0.0002004'
$ws.Range("E51").Value = '0.12%'
